# SALES ver 0.8 at 08/09/2022
# Update the tyre sales data table (columns H:J) on the data sheet and
# drop the trailing duplicated block of rows (97:134) that is no longer
# part of the report.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Holidays 2019")

# --- Rows 2-20: БНХ РОС block -> Sales value 4, Date 2022-09-20 ---------
$ws.Range("H2:H20").Value = 4
$ws.Range("I2:I20").Value = 44824

# --- Rows 21-39: БНХ Укр block -> Sales value 3, Date 2022-09-20, ------
# --- Contragent label recapitalised to "БНХ УКР" -----------------------
$ws.Range("H21:H39").Value = 3
$ws.Range("I21:I39").Value = 44824
$ws.Range("J21:J39").Value = "БНХ УКР"

# --- Rows 40-58: БНХ Польска block -> Sales value 5, Date 2022-09-20 ---
$ws.Range("H40:H58").Value = 5
$ws.Range("I40:I58").Value = 44824

# --- Rows 59-77: БНХ Польска block -> Sales value 6, Date 2022-07-27 ---
$ws.Range("H59:H77").Value = 6
$ws.Range("I59:I77").Value = 44769

# --- Rows 78-96: БНХ Польска block -> Sales value 12, Date 2022-05-15 --
$ws.Range("H78:H96").Value = 12
$ws.Range("I78:I96").Value = 44696

# --- Drop the trailing duplicated rows 97:134 ---------------------------
$ws.Range("A97:J134").EntireRow.Delete() | Out-Null
